# The workbook's 16 sheets each hold a cached statsmodels OLS summary
# (backward-elimination run) in cell B2. The model was re-run, which only
# changed the "Date:" / "Time:" stamp baked into that summary text -- all
# the numeric results stay identical. Replicate that here.

$wb = $excel.ActiveWorkbook

$oldDate = "Sun, 05 Jan 2020"
$newDate = "Wed, 08 Jan 2020"

$sheetCount = $wb.Worksheets.Count

for ($i = 1; $i -le $sheetCount; $i++) {
    $ws = $wb.Worksheets.Item($i)
    $cell = $ws.Range("B2")
    $summary = $cell.Text

    if ($summary -eq $null -or $summary -eq "") {
        continue
    }

    $summary = $summary.Replace($oldDate, $newDate)
    $summary = $summary.Replace("21:22:27", "19:07:32")
    $summary = $summary.Replace("21:22:28", "19:07:32")

    # The very last sheet's run landed one second later than the rest.
    if ($i -eq $sheetCount) {
        $summary = $summary.Replace("19:07:32", "19:07:33")
    }

    $cell.Value = $summary
}
